$wb = $excel.ActiveWorkbook

# --- Rename the three "Include from iso3166-1-*" sheets ---
$wb.Worksheets.Item(2).Name = "Include ValueSets"
$wb.Worksheets.Item(3).Name = "Include ValueSets 2"
$wb.Worksheets.Item(4).Name = "Include ValueSets 3"

# --- Metadata sheet updates ---
$ws1 = $wb.Worksheets.Item(1)

# Version 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date updated
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 "Contact" / "No display for ContactDetail" becomes
# "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# it is removed entirely (rows below shift up).
$ws1.Rows.Item(11).Delete()

# --- Include sheets: collapse the 4-row "Codes / All codes / (blank) /
#     System URI" block into a 2-row "ValueSet URL" block ---
foreach ($idx in 2, 3, 4) {
    $ws = $wb.Worksheets.Item($idx)
    $url = $ws.Range("B4").Text
    $ws.Rows.Item(4).Delete()
    $ws.Rows.Item(3).Delete()
    $ws.Range("A1").Value = "ValueSet URL"
    $ws.Range("A2").Value = $url
}
